$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume updates (Price column stored as text; force
# text format where the new value would otherwise be auto-parsed as a number)

$ws.Range("D2").Value = "29.731.90"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").Value = "2.096.26"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.53"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5161"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4382"
$ws.Range("E8").Value = "  -3.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.03"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09219"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.83"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "2.101.87"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.250"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.771"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.54"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001151"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.011"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06656"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.213"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").Value = "29.759.96"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.50"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.321"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "2.347.18"
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.518"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.25"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.12"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.140"
$ws.Range("E31").Value = "  -6.53%  "
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.634"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.171"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.952"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.298"
$ws.Range("E36").Value = "  +2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.21"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02575"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7123"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06736"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.43"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.321"
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7013"
$ws.Range("E44").Value = "  +7.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.28"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.317"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.624"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000355"
$ws.Range("E49").Value = "  -4.62%  "
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.05"
$ws.Range("E51").Value = "  -2.02%  "
